$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 35, shifting existing rows 35-43 down to 37-45.
$ws.Range("A35:R36").EntireRow.Insert()

# Row 35 (new data)
$ws.Range("A35").Value = 9
$ws.Range("B35").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C35").Value = "Metropolitana"
$ws.Range("D35").Value = 44452
$ws.Range("D35").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E35").Value = 13
$ws.Range("F35").Value = 100114002
$ws.Range("G35").Value = "Camote"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 970
$ws.Range("K35").Value = 14000
$ws.Range("L35").Value = 15000
$ws.Range("M35").Value = 14495
$ws.Range("N35").Value = "$/malla 18 kilos"
$ws.Range("O35").Value = "Perú"
$ws.Range("P35").Value = 805
$ws.Range("Q35").Value = 18
$ws.Range("R35").Value = "Hortaliza"

# Row 36 (new data)
$ws.Range("A36").Value = 9
$ws.Range("B36").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C36").Value = "Metropolitana"
$ws.Range("D36").Value = 44452
$ws.Range("D36").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E36").Value = 13
$ws.Range("F36").Value = 100114002
$ws.Range("G36").Value = "Camote"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Segunda"
$ws.Range("J36").Value = 340
$ws.Range("K36").Value = 11000
$ws.Range("L36").Value = 12000
$ws.Range("M36").Value = 11500
$ws.Range("N36").Value = "$/malla 18 kilos"
$ws.Range("O36").Value = "Perú"
$ws.Range("P36").Value = 639
$ws.Range("Q36").Value = 18
$ws.Range("R36").Value = "Hortaliza"
